$wb = $excel.ActiveWorkbook

# Plan2: add column F with new values C05 / teste1 / teste2, mapped to rows 2,3,4
$ws2 = $wb.Worksheets.Item("Plan2")
$ws2.Range("F3").Value = "teste1"
$ws2.Range("F4").Value = "teste2"
$ws2.Range("F2").Value = "C05"
$ws2.Range("F3").Select()

# Plan1: clear K5 then re-assert value (forces sst reinsertion) then move selection to K7
$ws1 = $wb.Worksheets.Item("Plan1")
$ws1.Range("K5").ClearContents()
$ws1.Range("K5").Value = "teste 1;teste 2"
$ws1.Range("K7").Select()
